$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033439050578729
$ws.Range("D2").Value = 1.037079819858861
$ws.Range("E2").Value = 1.054028787432045
$ws.Range("F2").Value = 1.059975294795178
$ws.Range("I2").Value = 1.039819553852152
$ws.Range("J2").Value = 1.038563771810583
$ws.Range("K2").Value = 1.03987172730686
$ws.Range("L2").Value = 1.056773126186569
$ws.Range("M2").Value = 1.062703334657551
$ws.Range("N2").Value = 1.040038651491483
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034258625655614
$ws.Range("D3").Value = 1.037686267510615
$ws.Range("E3").Value = 1.055080277588298
$ws.Range("F3").Value = 1.061043763695929
$ws.Range("I3").Value = 1.040037924140407
$ws.Range("J3").Value = 1.039026719630982
$ws.Range("K3").Value = 1.040288447557296
$ws.Range("L3").Value = 1.05763716867021
$ws.Range("M3").Value = 1.06358549608636
$ws.Range("N3").Value = 1.040502256750889
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034789197301381
$ws.Range("D4").Value = 1.038078829246122
$ws.Range("E4").Value = 1.055761813243849
$ws.Range("F4").Value = 1.061736101811262
$ws.Range("I4").Value = 1.040178007024959
$ws.Range("J4").Value = 1.039325862360939
$ws.Range("K4").Value = 1.040557548140531
$ws.Range("L4").Value = 1.058196808580137
$ws.Range("M4").Value = 1.064156678659055
$ws.Range("N4").Value = 1.040801824297835
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.035012308537722
$ws.Range("D5").Value = 1.038243896491074
$ws.Range("E5").Value = 1.056048605574049
$ws.Range("F5").Value = 1.062027391388771
$ws.Range("I5").Value = 1.040236605783694
$ws.Range("J5").Value = 1.039451521491132
$ws.Range("K5").Value = 1.040670546454864
$ws.Range("L5").Value = 1.058432211288108
$ws.Range("M5").Value = 1.06439689034692
$ws.Range("N5").Value = 1.040927661878406
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.0350497733115
$ws.Range("D6").Value = 1.038271613974671
$ws.Range("E6").Value = 1.056096775354225
$ws.Range("F6").Value = 1.062076313687928
$ws.Range("I6").Value = 1.040246427625807
$ws.Range("J6").Value = 1.039472614290835
$ws.Range("K6").Value = 1.04068951162375
$ws.Range("L6").Value = 1.058471744044615
$ws.Range("M6").Value = 1.064437227995208
$ws.Range("N6").Value = 1.040948784632303
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034792178294836
$ws.Range("D7").Value = 1.038081034749891
$ws.Range("E7").Value = 1.055765644298527
$ws.Range("F7").Value = 1.061739993131783
$ws.Range("I7").Value = 1.040178791172967
$ws.Range("J7").Value = 1.039327541820445
$ws.Range("K7").Value = 1.040559058547813
$ws.Range("L7").Value = 1.058199953532444
$ws.Range("M7").Value = 1.064159888040459
$ws.Range("N7").Value = 1.040803506142366
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033715976866858
$ws.Range("D8").Value = 1.037284739808287
$ws.Range("E8").Value = 1.054383904567084
$ws.Range("F8").Value = 1.060336188013781
$ws.Range("I8").Value = 1.039893604763288
$ws.Range("J8").Value = 1.038720312612016
$ws.Range("K8").Value = 1.040012672153113
$ws.Range("L8").Value = 1.057065019733272
$ws.Range("M8").Value = 1.0630013891071
$ws.Range("N8").Value = 1.04019541459881
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031821560594876
$ws.Range("D9").Value = 1.035882778440686
$ws.Range("E9").Value = 1.051957957654542
$ws.Range("F9").Value = 1.057869953747254
$ws.Range("I9").Value = 1.039381776135489
$ws.Range("J9").Value = 1.037647159735135
$ws.Range("K9").Value = 1.039045738657603
$ws.Range("L9").Value = 1.055069342436035
$ws.Range("M9").Value = 1.060962801431548
$ws.Range("N9").Value = 1.039120737721752
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030560033100994
$ws.Range("D10").Value = 1.034949041995969
$ws.Range("E10").Value = 1.050346663354076
$ws.Range("F10").Value = 1.056230860761687
$ws.Range("I10").Value = 1.039034347592983
$ws.Range("J10").Value = 1.036929674633136
$ws.Range("K10").Value = 1.038398399726624
$ws.Range("L10").Value = 1.053741774112775
$ws.Range("M10").Value = 1.059605700659386
$ws.Range("N10").Value = 1.038402233708611
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030014130393678
$ws.Range("D11").Value = 1.034544957121131
$ws.Range("E11").Value = 1.049650389214767
$ws.Range("F11").Value = 1.055522327335347
$ws.Range("I11").Value = 1.038882442519408
$ws.Range("J11").Value = 1.036618521870262
$ws.Range("K11").Value = 1.038117463157328
$ws.Range("L11").Value = 1.053167615370429
$ws.Range("M11").Value = 1.059018535787717
$ws.Range("N11").Value = 1.038090639073124
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029811411263334
$ws.Range("D12").Value = 1.034394897940769
$ws.Range("E12").Value = 1.049391977140153
$ws.Range("F12").Value = 1.055259328291106
$ws.Range("I12").Value = 1.038825798551825
$ws.Range("J12").Value = 1.036502875193636
$ws.Range("K12").Value = 1.038013016472435
$ws.Range("L12").Value = 1.052954451165561
$ws.Range("M12").Value = 1.058800507816873
$ws.Range("N12").Value = 1.037974828164952
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02985489276563
$ws.Range("D13").Value = 1.034427084494142
$ws.Range("E13").Value = 1.049447397659403
$ws.Range("F13").Value = 1.055315734237804
$ws.Range("I13").Value = 1.038837958805195
$ws.Range("J13").Value = 1.036527684981415
$ws.Range("K13").Value = 1.03803542489371
$ws.Range("L13").Value = 1.05300017089758
$ws.Range("M13").Value = 1.058847272321828
$ws.Range("N13").Value = 1.037999673185476
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029997372469381
$ws.Range("D14").Value = 1.034532552440742
$ws.Range("E14").Value = 1.049629024398793
$ws.Range("F14").Value = 1.055500584038442
$ws.Range("I14").Value = 1.038877764788757
$ws.Range("J14").Value = 1.036608963916243
$ws.Range("K14").Value = 1.038108831480395
$ws.Range("L14").Value = 1.053149993026897
$ws.Range("M14").Value = 1.059000512070483
$ws.Range("N14").Value = 1.038081067545715
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030085166037867
$ws.Range("D15").Value = 1.03459753953006
$ws.Range("E15").Value = 1.049740959158747
$ws.Range("F15").Value = 1.055614500219051
$ws.Range("I15").Value = 1.038902261480626
$ws.Range("J15").Value = 1.036659033213035
$ws.Range("K15").Value = 1.038154047230797
$ws.Range("L15").Value = 1.053242317176116
$ws.Range("M15").Value = 1.059094937596285
$ws.Range("N15").Value = 1.03813120794665
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030596270288542
$ws.Range("D16").Value = 1.03497586473507
$ws.Range("E16").Value = 1.05039290290424
$ws.Range("F16").Value = 1.056277909243248
$ws.Range("I16").Value = 1.039044398194022
$ws.Range("J16").Value = 1.036950314851024
$ws.Range("K16").Value = 1.038417031277675
$ws.Range("L16").Value = 1.053779893672538
$ws.Range("M16").Value = 1.059644678800493
$ws.Range("N16").Value = 1.038422903237976
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030916966310626
$ws.Range("D17").Value = 1.03521324055866
$ws.Range("E17").Value = 1.050802232556905
$ws.Range("F17").Value = 1.056694371302583
$ws.Range("I17").Value = 1.039133164689783
$ws.Range("J17").Value = 1.037132901114095
$ws.Range("K17").Value = 1.038581825042843
$ws.Range("L17").Value = 1.054117285646355
$ws.Range("M17").Value = 1.059989642976642
$ws.Range("N17").Value = 1.038605748794484
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.031104056269597
$ws.Range("D18").Value = 1.035351719871148
$ws.Range("E18").Value = 1.051041125244146
$ws.Range("F18").Value = 1.056937402751194
$ws.Range("I18").Value = 1.039184799173915
$ws.Range("J18").Value = 1.037239354505361
$ws.Range("K18").Value = 1.03867788516715
$ws.Range("L18").Value = 1.054314146925484
$ws.Range("M18").Value = 1.060190899915732
$ws.Range("N18").Value = 1.038712353361776
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031167854735449
$ws.Range("D19").Value = 1.035398941425653
$ws.Range("E19").Value = 1.051122604797617
$ws.Range("F19").Value = 1.057020289874803
$ws.Range("I19").Value = 1.039202381176967
$ws.Range("J19").Value = 1.037275644514
$ws.Range("K19").Value = 1.038710628760752
$ws.Range("L19").Value = 1.054381282750477
$ws.Range("M19").Value = 1.060259530969037
$ws.Range("N19").Value = 1.038748694906389
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030882555198841
$ws.Range("D20").Value = 1.035187770091751
$ws.Range("E20").Value = 1.050758301118158
$ws.Range("F20").Value = 1.056649676842351
$ws.Range("I20").Value = 1.039123655517503
$ws.Range("J20").Value = 1.037113316095978
$ws.Range("K20").Value = 1.03856415056381
$ws.Range("L20").Value = 1.054081079823161
$ws.Range("M20").Value = 1.059952626907232
$ws.Range("N20").Value = 1.038586135963394
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029955414262105
$ws.Range("D21").Value = 1.034501493755171
$ws.Range("E21").Value = 1.049575533905831
$ws.Range("F21").Value = 1.05544614535364
$ws.Range("I21").Value = 1.038866048977406
$ws.Range("J21").Value = 1.036585031239229
$ws.Range("K21").Value = 1.038087217666306
$ws.Range("L21").Value = 1.053105871269884
$ws.Range("M21").Value = 1.058955384811928
$ws.Range("N21").Value = 1.038057100881554
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029372793800085
$ws.Range("D22").Value = 1.034070213229286
$ws.Range("E22").Value = 1.048833126344908
$ws.Range("F22").Value = 1.054690489696622
$ws.Range("I22").Value = 1.038702810947575
$ws.Range("J22").Value = 1.036252469094445
$ws.Range("K22").Value = 1.037786805728157
$ws.Range("L22").Value = 1.052493320519318
$ws.Range("M22").Value = 1.05832879149547
$ws.Range("N22").Value = 1.037724066460379
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029681622124803
$ws.Range("D23").Value = 1.034298822917704
$ws.Range("E23").Value = 1.049226572359552
$ws.Range("F23").Value = 1.055090977027004
$ws.Range("I23").Value = 1.038789466742953
$ws.Range("J23").Value = 1.036428804974377
$ws.Range("K23").Value = 1.037946111135281
$ws.Range("L23").Value = 1.052817988009008
$ws.Range("M23").Value = 1.058660921116223
$ws.Range("N23").Value = 1.037900652757486
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030898103995991
$ws.Range("D24").Value = 1.035199279031549
$ws.Range("E24").Value = 1.050778151419102
$ws.Range("F24").Value = 1.056669871986144
$ws.Range("I24").Value = 1.039127952740442
$ws.Range("J24").Value = 1.037122165865137
$ws.Range("K24").Value = 1.038572137089579
$ws.Range("L24").Value = 1.054097439471522
$ws.Range("M24").Value = 1.059969352736117
$ws.Range("N24").Value = 1.03859499830024
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032311068767617
$ws.Range("D25").Value = 1.036245066183063
$ws.Range("E25").Value = 1.052584068581382
$ws.Range("F25").Value = 1.058506646003942
$ws.Range("I25").Value = 1.039515193431281
$ws.Range("J25").Value = 1.037924961137142
$ws.Range("K25").Value = 1.039296197406773
$ws.Range("L25").Value = 1.055584767996725
$ws.Range("M25").Value = 1.06148948392835
$ws.Range("N25").Value = 1.039398933633614
